$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H54").Value = 4250
$ws.Range("J54").Value = 5000
$ws.Range("L54").Value = 5000
$ws.Range("N54").Value = -5972
$ws.Range("H86").Value = 4092.5881
$ws.Range("I86").Value = 2422.8333
$ws.Range("K86").Value = 2422.8333
$ws.Range("M86").Value = -1299.8333
$ws.Range("H89").Value = 4092.5881
$ws.Range("I89").Value = 2422.8333
$ws.Range("K89").Value = 12114.1665
$ws.Range("M89").Value = -6498.166499999999
$ws.Range("H105").Value = 12000
$ws.Range("J105").Value = 12000
$ws.Range("L105").Value = 12000
$ws.Range("N105").Value = -18988
$ws.Range("H137").Value = 1813.15
$ws.Range("I137").Value = 1436.1333
$ws.Range("K137").Value = 4308.3999
$ws.Range("M137").Value = -1758.3999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 200.57143
$ws.Range("I5").Value = 184
$ws.Range("K5").Value = 184
$ws.Range("M5").Value = -72
$ws.Range("H32").Value = 2542.2253
$ws.Range("I32").Value = 2681.0833
$ws.Range("J32").Value = 1784.8182
$ws.Range("K32").Value = 2681.0833
$ws.Range("L32").Value = 1784.8182
$ws.Range("M32").Value = -2394.0833
$ws.Range("N32").Value = -2358.8182
$ws.Range("H122").Value = 1235.9166
$ws.Range("I122").Value = 1266.091
$ws.Range("J122").Value = 904
$ws.Range("K122").Value = 3798.273
$ws.Range("L122").Value = 2712
$ws.Range("M122").Value = -1348.273
$ws.Range("N122").Value = -7612
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H132").Value = 2187.5625
$ws.Range("I132").Value = 1928.3572
$ws.Range("K132").Value = 5785.071599999999
$ws.Range("M132").Value = -3255.071599999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 200.57143
$ws.Range("I4").Value = 184
$ws.Range("K4").Value = 184
$ws.Range("M4").Value = -69
$ws.Range("H22").Value = 197.8
$ws.Range("I22").Value = 197.8
$ws.Range("K22").Value = 197.8
$ws.Range("M22").Value = -24.80000000000001
$ws.Range("H64").Value = 400
$ws.Range("I64").Value = 433.33334
$ws.Range("J64").Value = 200
$ws.Range("K64").Value = 433.33334
$ws.Range("L64").Value = 200
$ws.Range("M64").Value = -208.33334
$ws.Range("N64").Value = -650
$ws.Range("H67").Value = 400
$ws.Range("I67").Value = 433.33334
$ws.Range("J67").Value = 200
$ws.Range("K67").Value = 433.33334
$ws.Range("L67").Value = 200
$ws.Range("M67").Value = 346.66666
$ws.Range("N67").Value = -1760
$ws.Range("H107").Value = 1775.7241
$ws.Range("I107").Value = 1579.7142
$ws.Range("J107").Value = 2290.25
$ws.Range("K107").Value = 1579.7142
$ws.Range("L107").Value = 2290.25
$ws.Range("M107").Value = 340.2858000000001
$ws.Range("N107").Value = -6130.25

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1468.2424
$ws.Range("J31").Value = 2444.75
$ws.Range("L31").Value = 2444.75
$ws.Range("N31").Value = -3034.75
$ws.Range("H34").Value = 1468.2424
$ws.Range("J34").Value = 2444.75
$ws.Range("L34").Value = 2444.75
$ws.Range("N34").Value = -2848.75
$ws.Range("H62").Value = 7410418.5
$ws.Range("I62").Value = 3056.5217
$ws.Range("J62").Value = 50002750
$ws.Range("K62").Value = 3056.5217
$ws.Range("L62").Value = 50002750
$ws.Range("M62").Value = -2432.5217
$ws.Range("N62").Value = -50003998
$ws.Range("H65").Value = 7410418.5
$ws.Range("I65").Value = 3056.5217
$ws.Range("J65").Value = 50002750
$ws.Range("K65").Value = 15282.6085
$ws.Range("L65").Value = 250013750
$ws.Range("M65").Value = -12162.6085
$ws.Range("N65").Value = -250019990
$ws.Range("H132").Value = 11517.818
$ws.Range("J132").Value = 2938
$ws.Range("L132").Value = 8814
$ws.Range("N132").Value = -13874
$ws.Range("H133").Value = 32826
$ws.Range("J133").Value = 32826
$ws.Range("L133").Value = 32826
$ws.Range("N133").Value = -37886

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 3599
$ws.Range("J55").Value = 3599
$ws.Range("L55").Value = 10797
$ws.Range("N55").Value = -11151
$ws.Range("H113").Value = 538.3333
$ws.Range("I113").Value = 457.69232
$ws.Range("J113").Value = 568.2857
$ws.Range("K113").Value = 1373.07696
$ws.Range("L113").Value = 1704.8571
$ws.Range("M113").Value = 796.9230400000001
$ws.Range("N113").Value = -6044.8571
$ws.Range("H122").Value = 908.6
$ws.Range("I122").Value = 636
$ws.Range("J122").Value = 1999
$ws.Range("K122").Value = 5724
$ws.Range("L122").Value = 17991
$ws.Range("M122").Value = -3274
$ws.Range("N122").Value = -22891
$ws.Range("H131").Value = 12049401
$ws.Range("J131").Value = 1257.4736
$ws.Range("L131").Value = 3772.4208
$ws.Range("N131").Value = -13852.4208
$ws.Range("H136").Value = 2543.2
$ws.Range("I136").Value = 1020
$ws.Range("J136").Value = 4066.4
$ws.Range("K136").Value = 3060
$ws.Range("L136").Value = 12199.2
$ws.Range("M136").Value = 2040
$ws.Range("N136").Value = -22399.2

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 56255800
$ws.Range("I70").Value = 125003000
$ws.Range("J70").Value = 33340068
$ws.Range("K70").Value = 125003000
$ws.Range("L70").Value = 33340068
$ws.Range("M70").Value = -125002730
$ws.Range("N70").Value = -33340608
$ws.Range("H73").Value = 56255800
$ws.Range("I73").Value = 125003000
$ws.Range("J73").Value = 33340068
$ws.Range("K73").Value = 125003000
$ws.Range("L73").Value = 33340068
$ws.Range("M73").Value = -125002064
$ws.Range("N73").Value = -33341940
$ws.Range("H102").Value = 1749.0454
$ws.Range("I102").Value = 1836.7858
$ws.Range("J102").Value = 1595.5
$ws.Range("K102").Value = 1836.7858
$ws.Range("L102").Value = 1595.5
$ws.Range("M102").Value = -214.7858000000001
$ws.Range("N102").Value = -4839.5
$ws.Range("H132").Value = 2796.9412
$ws.Range("I132").Value = 2569
$ws.Range("J132").Value = 4506.5
$ws.Range("K132").Value = 7707
$ws.Range("L132").Value = 13519.5
$ws.Range("M132").Value = -5177
$ws.Range("N132").Value = -18579.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 1000
$ws.Range("I32").Value = 1000
$ws.Range("K32").Value = 1000
$ws.Range("M32").Value = -683
$ws.Range("H46").Value = 5818.1875
$ws.Range("I46").Value = 577.8
$ws.Range("J46").Value = 8200.182000000001
$ws.Range("K46").Value = 577.8
$ws.Range("L46").Value = 8200.182000000001
$ws.Range("M46").Value = -389.8
$ws.Range("N46").Value = -8576.182000000001
$ws.Range("H55").Value = 772.63635
$ws.Range("I55").Value = 250.25
$ws.Range("K55").Value = 250.25
$ws.Range("M55").Value = -77.25
$ws.Range("H68").Value = 1561.909
$ws.Range("I68").Value = 1310.375
$ws.Range("K68").Value = 1310.375
$ws.Range("M68").Value = -561.375
$ws.Range("H71").Value = 1561.909
$ws.Range("I71").Value = 1310.375
$ws.Range("K71").Value = 6551.875
$ws.Range("M71").Value = -2807.875
$ws.Range("H122").Value = 23619388
$ws.Range("I122").Value = 23619388
$ws.Range("K122").Value = 70858164
$ws.Range("M122").Value = -70855714
$ws.Range("H132").Value = 35061.8
$ws.Range("I132").Value = 1397.6522
$ws.Range("K132").Value = 4192.9566
$ws.Range("M132").Value = -1662.9566
$ws.Range("H136").Value = 10134.546
$ws.Range("I136").Value = 10134.546
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 30403.638
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -27853.638
$ws.Range("N136").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 58502.8
$ws.Range("J123").Value = 58502.8
$ws.Range("L123").Value = 58502.8
$ws.Range("N123").Value = -68302.8
$ws.Range("H136").Value = 705.5333000000001
$ws.Range("I136").Value = 381.06668
$ws.Range("K136").Value = 1143.20004
$ws.Range("M136").Value = 1406.79996
